$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-7 (row1 header, row4 unchanged) reflecting
# the reordering described by the diff.
$ws.Range("A2").Value = "concept:name"
$ws.Range("B2").Value = "str"

$ws.Range("A3").Value = "SubProcessID"
$ws.Range("B3").Value = "str"

$ws.Range("A5").Value = "stream:datastream"
$ws.Range("B5").Value = "dict"

$ws.Range("A6").Value = "org:resource"
$ws.Range("B6").Value = "str"

$ws.Range("A7").Value = "time:timestamp"
$ws.Range("B7").Value = "datetime"
